# The underlying commit only changes one real input: the "Actual Words
# Tagged" entry for 2020-09-25 (row 4) in Table1 on the "tagging" sheet,
# which went from 4665 to 5259. Every other changed cell in the diff
# (Target/Actual Words Remaining, the Table2 rollups G2:J2, the mirrored
# Q2/Q5/Q8 cells on "graphs", and the two chart caches) is a formula or a
# chart series that is derived from that single cell, so touching it and
# letting Excel recalculate reproduces the rest of the diff automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tagging")

# Table1[[#This Row],[Actual Words Tagged]] for the 2020-09-25 row.
$ws.Range("C4").Value = 5259

# Make sure everything (formulas, table calculated columns, and the
# dependent charts on the "graphs" sheet) is fully up to date before save.
$wb.RefreshAll()
$excel.CalculateFullRebuild()
